$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Fixed k loss" row (row 16) values from 2.4 to 2.2 for all 4 datasets
$ws.Range("B16:E16").Value = 2.2

# Update the active selection to reflect the new cell (G19 instead of G20)
$ws.Range("G19").Select()
